# Weekly fruit/vegetable update: insert a new daily record for
# "Terminal Hortofrutícola Agro Chillán - Repollo" before the existing
# row 194, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 194 (pushes old rows 194-262 down to 195-263).
$ws.Range("A194").EntireRow.Insert()

# Populate the newly inserted row 194 with the new record's data.
$ws.Cells.Item(194, 1).Value = 7
$ws.Cells.Item(194, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(194, 3).Value = "Ñuble"
$ws.Cells.Item(194, 4).Value = 44795
$ws.Cells.Item(194, 5).Value = 16
$ws.Cells.Item(194, 6).Value = 100112006
$ws.Cells.Item(194, 7).Value = "Repollo"
$ws.Cells.Item(194, 8).Value = "Crespo record"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 160
$ws.Cells.Item(194, 11).Value = 1200
$ws.Cells.Item(194, 12).Value = 1300
$ws.Cells.Item(194, 13).Value = 1250
$ws.Cells.Item(194, 14).Value = "$/unidad"
$ws.Cells.Item(194, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(194, 16).Value = 1250
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"
